$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text formatting
# (so numeric-looking strings like "303.50" or "1.00" are not
# silently converted to numbers and lose trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '41.739.38'
$ws.Range('E2').Value = '  +1.29%  '

$ws.Range('D3').Value = '2.265.45'
$ws.Range('E3').Value = '  +0.74%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = '303.50'
$ws.Range('E5').Value = '  +0.32%  '

$ws.Range('D6').Value = '92.15'
$ws.Range('E6').Value = '  +1.25%  '

$ws.Range('E7').Value = '  +1.97%  '

$ws.Range('E8').Value = '  +0.01%  '

$ws.Range('D9').Value = '0.482'
$ws.Range('E9').Value = '  -0.07%  '

$ws.Range('D10').Value = '32.42'
$ws.Range('E10').Value = '  +1.97%  '

$ws.Range('D11').Value = '53.47'
$ws.Range('E11').Value = '  -1.29%  '

$ws.Range('D12').Value = '0.0796'
$ws.Range('E12').Value = '  +0.28%  '

$ws.Range('D13').Value = '0.113'
$ws.Range('E13').Value = '  -1.25%  '

$ws.Range('D14').Value = '6.64'
$ws.Range('E14').Value = '  +1.15%  '

$ws.Range('D15').Value = '2.615.97'
$ws.Range('E15').Value = '  +0.80%  '

$ws.Range('D16').Value = '14.24'
$ws.Range('E16').Value = '  +1.01%  '

$ws.Range('D17').Value = '2.255.66'
$ws.Range('E17').Value = '  -0.35%  '

$ws.Range('D18').Value = '0.769'
$ws.Range('E18').Value = '  +2.64%  '

$ws.Range('D19').Value = '41.647.88'
$ws.Range('E19').Value = '  +1.27%  '

$ws.Range('D20').Value = '12.41'
$ws.Range('E20').Value = '  +4.51%  '

$ws.Range('E21').Value = '  +0.32%  '

$ws.Range('D22').Value = '5.95'
$ws.Range('E22').Value = '  +1.69%  '

$ws.Range('D23').Value = '67.12'
$ws.Range('E23').Value = '  +0.38%  '

$ws.Range('D24').Value = '239.66'
$ws.Range('E24').Value = '  -0.27%  '

$ws.Range('D25').Value = '2.60'
$ws.Range('E25').Value = '  +1.77%  '

$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').Value = '1.93'
$ws.Range('E26').Value = '  +4.34%  '

$ws.Range('B27').Value = 'Dai'
$ws.Range('C27').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D27').Value = '1.00'
$ws.Range('E27').Value = '  -0.08%  '

$ws.Range('D28').Value = '23.94'
$ws.Range('E28').Value = '  +0.90%  '

$ws.Range('E29').Value = '  -0.38%  '

$ws.Range('D30').Value = '2.12'
$ws.Range('E30').Value = '  -2.86%  '

$ws.Range('D31').Value = '35.35'
$ws.Range('E31').Value = '  +6.93%  '

$ws.Range('D32').Value = '159.95'
$ws.Range('E32').Value = '  +0.65%  '

$ws.Range('D33').Value = '5.24'
$ws.Range('E33').Value = '  +1.49%  '

$ws.Range('D34').Value = '1.00'
$ws.Range('E34').Value = '  -0.01%  '

$ws.Range('D36').Value = '3.01'
$ws.Range('E36').Value = '  -0.60%  '

$ws.Range('D37').Value = '16.93'
$ws.Range('E37').Value = '  +2.40%  '

$ws.Range('E38').Value = '  +0.43%  '

$ws.Range('E39').Value = '  +1.43%  '

$ws.Range('E40').Value = '  +1.01%  '

$ws.Range('E41').Value = '  +0.63%  '

$ws.Range('E42').Value = '  -0.12%  '

$ws.Range('D43').Value = '2.017.10'
$ws.Range('E43').Value = '  -2.83%  '

$ws.Range('D44').Value = '19.23'
$ws.Range('E44').Value = '  -4.11%  '

$ws.Range('E45').Value = '  +1.25%  '

$ws.Range('E46').Value = '  +0.86%  '

$ws.Range('E47').Value = '  +5.53%  '

$ws.Range('D48').Value = '2.88'
$ws.Range('E48').Value = '  -1.93%  '

$ws.Range('E49').Value = '  +2.23%  '

$ws.Range('E50').Value = '  +0.72%  '

$ws.Range('D51').Value = '52.21'
$ws.Range('E51').Value = '  +2.83%  '
